$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Delete the "Phụ cấp tại SÓC TRĂNG" row (row 24) first, then the
# "Phụ cấp tại CẦN THƠ" row (row 3), so row indices of the first
# deletion are not disturbed by the second.
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(3).Delete()
